$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the Fecha (D) and Volumen (J) values among rows 2, 3 and 5:
#   row2 <- old row3, row3 <- old row5, row5 <- old row2

$d2 = $ws.Range("D2").Value2
$d3 = $ws.Range("D3").Value2
$d5 = $ws.Range("D5").Value2

$j2 = $ws.Range("J2").Value2
$j3 = $ws.Range("J3").Value2
$j5 = $ws.Range("J5").Value2

$ws.Range("D2").Value2 = $d3
$ws.Range("D3").Value2 = $d5
$ws.Range("D5").Value2 = $d2

$ws.Range("J2").Value2 = $j3
$ws.Range("J3").Value2 = $j5
$ws.Range("J5").Value2 = $j2
